$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.015.73"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "3.880.07"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'599.13"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'171.56"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("D7").Value = "3.879.76"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'37.04"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "4.532.73"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "3.887.38"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "68.184.70"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "'18.10"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").Value = "'7.36"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'467.32"
$ws.Range("E22").Value = "  -5.50%  "
$ws.Range("D23").Value = "'0.739"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").Value = "'83.31"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "'12.08"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "'2.95"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "4.033.44"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "'9.49"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "3.852.28"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "'3.81"
$ws.Range("E37").Value = "  +15.66%  "
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "'1.03"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D41").Value = "'5.90"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "'0.000301"
$ws.Range("E44").Value = "  +11.51%  "
$ws.Range("D45").Value = "'425.64"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'8.63"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'47.23"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'27.50"
$ws.Range("E50").Value = "  +6.54%  "
$ws.Range("D51").Value = "'143.32"
$ws.Range("E51").Value = "  +0.29%  "
